$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1) updates - column F (想去人数)
$ws1.Range("F2").Value = 3134
$ws1.Range("F3").Value = 530
$ws1.Range("F4").Value = 1092
$ws1.Range("F6").Value = 35
$ws1.Range("F8").Value = 36
$ws1.Range("F9").Value = 1125
$ws1.Range("F10").Value = 15661
$ws1.Range("F11").Value = 239
$ws1.Range("F12").Value = 175
$ws1.Range("F13").Value = 1024
$ws1.Range("F14").Value = 6167
$ws1.Range("F15").Value = 621
$ws1.Range("F17").Value = 65
$ws1.Range("F19").Value = 113
$ws1.Range("F23").Value = 13
$ws1.Range("F26").Value = 209
$ws1.Range("F29").Value = 4996
$ws1.Range("F31").Value = 11046
$ws1.Range("F32").Value = 1228
$ws1.Range("F35").Value = 165
$ws1.Range("F36").Value = 3802
$ws1.Range("F37").Value = 264
$ws1.Range("F38").Value = 73

# Sheet "演出" (sheet2) update
$ws2.Range("F3").Value = 20

# Sheet "全部类型" (sheet4) updates
$ws4.Range("F3").Value = 3134
$ws4.Range("F4").Value = 530
$ws4.Range("F5").Value = 1092
$ws4.Range("F7").Value = 35
$ws4.Range("F9").Value = 36
$ws4.Range("F10").Value = 1125
$ws4.Range("F11").Value = 15662
$ws4.Range("F12").Value = 239
$ws4.Range("F13").Value = 175
$ws4.Range("F14").Value = 1024
$ws4.Range("F15").Value = 6167
$ws4.Range("F16").Value = 621
$ws4.Range("F18").Value = 65
$ws4.Range("F20").Value = 113
$ws4.Range("F24").Value = 13
$ws4.Range("F27").Value = 209
$ws4.Range("F30").Value = 4996
$ws4.Range("F32").Value = 20
$ws4.Range("F33").Value = 11046
$ws4.Range("F34").Value = 1228
$ws4.Range("F37").Value = 165
$ws4.Range("F38").Value = 3802
$ws4.Range("F39").Value = 264
$ws4.Range("F40").Value = 73
